# Actualización automática 2025-07-01 08:30:08
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": zero-out a handful of client figures and refresh
# the "x de 27" progress labels on row 29 to reflect the reset counts.
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("D13").Value = 0
$wsGrupo.Range("E13").Value = 0

$wsGrupo.Range("E14").Value = 0
$wsGrupo.Range("N14").Value = 0
$wsGrupo.Range("Q14").Value = 0

$wsGrupo.Range("C17").Value = 0
$wsGrupo.Range("D17").Value = 0
$wsGrupo.Range("E17").Value = 0
$wsGrupo.Range("G17").Value = 0
$wsGrupo.Range("H17").Value = 0
$wsGrupo.Range("I17").Value = 0
$wsGrupo.Range("L17").Value = 0
$wsGrupo.Range("M17").Value = 0

$wsGrupo.Range("E18").Value = 0

$wsGrupo.Range("D19").Value = 0
$wsGrupo.Range("M19").Value = 0

$wsGrupo.Range("L21").Value = 0
$wsGrupo.Range("O21").Value = 0

$wsGrupo.Range("C29").Value = "0 de 27"
$wsGrupo.Range("D29").Value = "0 de 27"
$wsGrupo.Range("E29").Value = "0 de 27"
$wsGrupo.Range("G29").Value = "0 de 27"
$wsGrupo.Range("H29").Value = "0 de 27"
$wsGrupo.Range("I29").Value = "0 de 27"
$wsGrupo.Range("L29").Value = "0 de 27"
$wsGrupo.Range("M29").Value = "0 de 27"
$wsGrupo.Range("N29").Value = "0 de 27"
$wsGrupo.Range("O29").Value = "0 de 27"
$wsGrupo.Range("Q29").Value = "0 de 27"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL": roll the monthly columns forward by one month
# (marzo/abril/mayo/junio -> abril/mayo/junio/julio), which shifts every
# client's figures left by one column (C<-D, D<-E, E<-F, F<-0) and widens /
# narrows the columns to match the new header widths.
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("C1").Value = "abril"
$wsMensual.Range("D1").Value = "mayo"
$wsMensual.Range("E1").Value = "junio"
$wsMensual.Range("F1").Value = "julio"

# Column widths: C and D swap (13 <-> 11), E is unchanged (13), F matches
# the old D width (11). ColumnWidth on this host adds 5/6 of a character to
# whatever is assigned, so subtract that back off to land on exact values.
$wsMensual.Columns(3).ColumnWidth = 11 - 0.8333333333333333
$wsMensual.Columns(4).ColumnWidth = 13 - 0.8333333333333333
$wsMensual.Columns(6).ColumnWidth = 11 - 0.8333333333333333

$rows = 2..29
foreach ($r in $rows) {
    $d = $wsMensual.Cells.Item($r, 4).Value2
    $e = $wsMensual.Cells.Item($r, 5).Value2
    $f = $wsMensual.Cells.Item($r, 6).Value2

    $wsMensual.Cells.Item($r, 3).Value = $d
    $wsMensual.Cells.Item($r, 4).Value = $e
    $wsMensual.Cells.Item($r, 5).Value = $f
    $wsMensual.Cells.Item($r, 6).Value = 0
}
